$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A168").Value = "IMX-USD"
$ws.Range("A169").Value = "TAO-USD"
$ws.Range("A170").Value = "GRT-USD"
